$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty row 12 with a new work-log entry:
# Who?, Task, Notes, Date, Time spent
$ws.Range("A12").Value = "Laurent"
$ws.Range("B12").Value = "Initial plan"
$ws.Range("C12").Value = "Worked on slides"
$ws.Range("D12").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("D12").Value = 42450
$ws.Range("E12").Value = "1h"

# Update the active selection to the newly-entered cell, matching the
# author's last cursor position after typing the row.
$ws.Range("E12").Select()
